$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Cells.Item($r, 58)
    $cell.NumberFormat = "@"
    $cell.Value = "2008-05-26"
    $cell.Style = "Normal"
}
